$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.231.38'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.904.13'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('D4').Value = "'0.9998"
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "'306.03"
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').Value = "'0.9993"
$ws.Range('D7').Value = "'0.5394"
$ws.Range('E7').Value = '  +3.42%  '
$ws.Range('D8').Value = "'0.3808"
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('D11').Value = "'0.9044"
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').Value = "'0.08180"
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = "'95.71"
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').Value = "'5.347"
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = "'0.9998"
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = "'0.000008660"
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').Value = "'0.9992"
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').Value = '27.249.84'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('B20').Value = 'WrappedEther'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D20').Value = '1.147.72'
$ws.Range('E20').Value = '  -39.48%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'5.047"
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').Value = "'10.82"
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('D24').Value = "'148.45"
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('D25').Value = "'2.306"
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('D27').Value = "'1.753"
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').Value = "'116.73"
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('D29').Value = "'4.858"
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('D30').Value = "'4.700"
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('D31').Value = "'0.09214"
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = "'0.8279"
$ws.Range('E32').Value = '  +4.97%  '
$ws.Range('D33').Value = "'0.05081"
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('D34').Value = "'1.219"
$ws.Range('E34').Value = '  +0.50%  '
$ws.Range('D35').Value = "'3.004"
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('D36').Value = "'3.322"
$ws.Range('E36').Value = '  -3.02%  '
$ws.Range('D37').Value = "'2.693"
$ws.Range('E37').Value = '  +4.38%  '
$ws.Range('D38').Value = "'0.5942"
$ws.Range('E38').Value = '  +4.39%  '
$ws.Range('D39').Value = "'0.02000"
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('D41').Value = "'9.278"
$ws.Range('E41').Value = '  +2.97%  '
$ws.Range('D42').Value = "'6.665"
$ws.Range('E42').Value = '  +1.71%  '
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').Value = "'0.5115"
$ws.Range('E44').Value = '  +5.39%  '
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('D46').Value = "'10.23"
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('D47').Value = "'0.9988"
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').Value = "'1.643"
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').Value = "'38.25"
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('D50').Value = "'0.06107"
$ws.Range('E50').Value = '  +2.97%  '
$ws.Range('D51').Value = "'63.50"
$ws.Range('E51').Value = '  +0.16%  '
